$wb = $excel.ActiveWorkbook

# --- 1. Rename "Reconcile New Statement" -> "Reconciliation" ---
$wsRecon = $wb.Worksheets.Item("Reconcile New Statement")
$wsRecon.Name = "Reconciliation"

# Change its selection from L25 to F22
$wsRecon.Range("F22").Select()

# --- 2. InvoiceNonPO: delete column G (ItemNo), becomes the active sheet, ---
#        topLeftCell scrolled to E1, new selection I15 ---
$wsInvNonPO = $wb.Worksheets.Item("InvoiceNonPO")
$wsInvNonPO.Columns("G:G").Delete()

# Make InvoiceNonPO the active sheet (this also clears tabSelected on Login
# and sets tabSelected + activeTab on InvoiceNonPO automatically)
$wsInvNonPO.Activate()
$wsInvNonPO.Range("I15").Select()

Write-Output "done"
